$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table")

$ws.Range("A1").Value = "Sample"
$ws.Range("B1").Value = "Year"
$ws.Range("C1").Value = "Mean"
$ws.Range("D1").Value = "Median"
$ws.Range("E1").Value = "<2500gm (%)"
$ws.Range("F1").Value = "Maternal Age (mean)"
$ws.Range("G1").Value = "Source"
$ws.Range("H1").Value = "Page #"

$ws.Range("A2").Value = "Antebellum US South (estimated enslaved)"
$ws.Range("B2").Value = "<1850"
$ws.Range("C2").Value = 2320
$ws.Range("G2").Value = "Steckel, 1986"
$ws.Range("H2").Value = "174, 182"

$ws.Range("A3").Value = "Rio de Janeiro (Black singletons)"
$ws.Range("B3").Value = "1922-26"
$ws.Range("C3").Value = 3037
$ws.Range("G3").Value = "Laranjeiras"

$ws.Range("A4").Value = "Rio de Janeiro (female singletons)"
$ws.Range("B4").Value = "1922-26"
$ws.Range("C4").Value = 3038
$ws.Range("G4").Value = "Laranjeiras"

$ws.Range("A5").Value = "Rio de Janeiro (Mixed-race singletons)"
$ws.Range("B5").Value = "1922-26"
$ws.Range("C5").Value = 3064
$ws.Range("G5").Value = "Laranjeiras"

$ws.Range("A6").Value = "Rio de Janeiro (singletons all)"
$ws.Range("B6").Value = "1922-26"
$ws.Range("C6").Value = 3087
$ws.Range("F6").Value = 25.3
$ws.Range("G6").Value = "Laranjeiras"

$ws.Range("A7").Value = "Riberão Preto, São Paulo, Brazil"
$ws.Range("B7").Value = 1994
$ws.Range("C7").Value = 3115
$ws.Range("D7").Value = 3150
$ws.Range("G7").Value = "Silva, 1998"
$ws.Range("H7").Value = 77

$ws.Range("A8").Value = "Rio de Janeiro (Black singletons, mothers <K-12)"
$ws.Range("B8").Value = "1999-2001"
$ws.Range("C8").Value = 3122
$ws.Range("G8").Value = "Leal, 2006"
$ws.Range("H8").Value = 469

$ws.Range("A9").Value = "Rio de Janeiro (White singletons)"
$ws.Range("B9").Value = "1922-26"
$ws.Range("C9").Value = 3133
$ws.Range("G9").Value = "Laranjeiras"

$ws.Range("A10").Value = "Rio de Janeiro (male singletons)"
$ws.Range("B10").Value = "1922-26"
$ws.Range("C10").Value = 3137
$ws.Range("G10").Value = "Laranjeiras"

$ws.Range("A11").Value = "Rio de Janeiro (Mixed-race singletons, mothers <K-12)"
$ws.Range("B11").Value = "1999-2001"
$ws.Range("C11").Value = 3154
$ws.Range("G11").Value = "Leal, 2006"
$ws.Range("H11").Value = 469

$ws.Range("A12").Value = "São Paulo, Brazil (live)"
$ws.Range("B12").Value = "1993-98"
$ws.Range("C12").Value = 3155
$ws.Range("E12").Value = 9
$ws.Range("G12").Value = "Monteiro, 2000"
$ws.Range("H12").Value = 31

$ws.Range("A13").Value = "Pelotas, Rio Grande do Sul, Brazil (live singletons)"
$ws.Range("B13").Value = 2004
$ws.Range("C13").Value = 3167
$ws.Range("E13").Value = 9
$ws.Range("G13").Value = "Silveira, 2019"
$ws.Range("H13").Value = "i48"

$ws.Range("A14").Value = "Pelotas, Rio Grande do Sul, Brazil (live singletons)"
$ws.Range("B14").Value = 1993
$ws.Range("C14").Value = 3169
$ws.Range("E14").Value = 9.1
$ws.Range("G14").Value = "Silveira, 2019"
$ws.Range("H14").Value = "i48"

$ws.Range("A15").Value = "Baltimore (Black singletons)"
$ws.Range("B15").Value = "1897-1935"
$ws.Range("C15").Value = 3183
$ws.Range("D15").Value = 3175
$ws.Range("E15").Value = 11.4
$ws.Range("G15").Value = "Costa, 2004"
$ws.Range("H15").Value = 1065

$ws.Range("A16").Value = "Rio de Janeiro (Black singletons, mothers  >=K-12)"
$ws.Range("B16").Value = "1999-2001"
$ws.Range("C16").Value = 3185
$ws.Range("G16").Value = "Leal, 2006"
$ws.Range("H16").Value = 470

$ws.Range("A17").Value = "Rio de Janeiro (White singletons, mothers <K-12)"
$ws.Range("B17").Value = "1999-2001"
$ws.Range("C17").Value = 3186
$ws.Range("G17").Value = "Leal, 2006"
$ws.Range("H17").Value = 469

$ws.Range("A18").Value = "Pelotas, Rio Grande do Sul, Brazil (live singletons)"
$ws.Range("B18").Value = 2015
$ws.Range("C18").Value = 3198
$ws.Range("E18").Value = 8.3
$ws.Range("G18").Value = "Silveira, 2019"
$ws.Range("H18").Value = "i48"

$ws.Range("A19").Value = "Pelotas, Rio Grande do Sul, Brazil (live singletons)"
$ws.Range("B19").Value = 1982
$ws.Range("C19").Value = 3201
$ws.Range("E19").Value = 8.2
$ws.Range("G19").Value = "Silveira, 2019"
$ws.Range("H19").Value = "i48"

$ws.Range("A20").Value = "Rio de Janeiro (Mixed-race singletons, mothers >=K-12)"
$ws.Range("B20").Value = "1999-2001"
$ws.Range("C20").Value = 3210
$ws.Range("G20").Value = "Leal, 2006"
$ws.Range("H20").Value = 470

$ws.Range("A21").Value = "Rio de Janeiro (White singletons, mothers >=K-12)"
$ws.Range("B21").Value = "1999-2001"
$ws.Range("C21").Value = 3218
$ws.Range("G21").Value = "Leal, 2006"
$ws.Range("H21").Value = 470

$ws.Range("A22").Value = "Riberão Preto, São Paulo, Brazil"
$ws.Range("B22").Value = "1978-79"
$ws.Range("C22").Value = 3234
$ws.Range("D22").Value = 3250
$ws.Range("G22").Value = "Silva, 1998"
$ws.Range("H22").Value = 77

$ws.Range("A23").Value = "Boston (in hospital)"
$ws.Range("B23").Value = "1886-1900"
$ws.Range("C23").Value = 3330
$ws.Range("E23").Value = 6.9
$ws.Range("G23").Value = "Ward, 1993"
$ws.Range("H23").Value = "148-9"

$ws.Range("A24").Value = "Philadelphia (all)"
$ws.Range("B24").Value = "1848-73"
$ws.Range("C24").Value = 3375
$ws.Range("D24").Value = 3453
$ws.Range("E24").Value = 9.6
$ws.Range("G24").Value = "Goldin, 1989"
$ws.Range("H24").Value = "363-5"

$ws.Range("A25").Value = "Philadelphia (live)"
$ws.Range("B25").Value = "1848-73"
$ws.Range("C25").Value = 3403
$ws.Range("D25").Value = 3461
$ws.Range("E25").Value = 8.1
$ws.Range("G25").Value = "Goldin, 1989"
$ws.Range("H25").Value = "363-5"

$ws.Range("A26").Value = "Wellington, NZ (singleton live female)"
$ws.Range("B26").Value = "1907-22"
$ws.Range("C26").Value = 3403
$ws.Range("G26").Value = "Roberts, 2014"
$ws.Range("H26").Value = "156, 158"

$ws.Range("A27").Value = "Baltimore (white singletons)"
$ws.Range("B27").Value = "1897-1935"
$ws.Range("C27").Value = 3423
$ws.Range("D27").Value = 3443
$ws.Range("E27").Value = 6
$ws.Range("G27").Value = "Costa, 2004"
$ws.Range("H27").Value = 1065

$ws.Range("A28").Value = "New York (singeltons)"
$ws.Range("B28").Value = "1910-31"
$ws.Range("C28").Value = 3463
$ws.Range("D28").Value = 3467
$ws.Range("E28").Value = 5.5
$ws.Range("F28").Value = 27
$ws.Range("G28").Value = "Costa, 1998"
$ws.Range("H28").Value = "991-2"

$ws.Range("A29").Value = "Wellington, NZ (singleton live)"
$ws.Range("B29").Value = "1907-22"
$ws.Range("C29").Value = 3467
$ws.Range("E29").Value = 4.2
$ws.Range("F29").Value = 27.75
$ws.Range("G29").Value = "Roberts, 2014"
$ws.Range("H29").Value = "156, 158"

$ws.Range("A30").Value = "Boston (at home)"
$ws.Range("B30").Value = "1884-1900"
$ws.Range("C30").Value = 3479
$ws.Range("E30").Value = 4.7
$ws.Range("G30").Value = "Ward, 1993"
$ws.Range("H30").Value = "148-9"

$ws.Range("A31").Value = "Boston"
$ws.Range("B31").Value = "1872-1900"
$ws.Range("C31").Value = 3480
$ws.Range("E31").Value = 6.5
$ws.Range("G31").Value = "Ward, 1993"
$ws.Range("H31").Value = "148-9"

$ws.Range("A32").Value = "Wellington, NZ (singleton live male)"
$ws.Range("B32").Value = "1907-22"
$ws.Range("C32").Value = 3531
$ws.Range("G32").Value = "Roberts, 2014"
$ws.Range("H32").Value = "156, 158"

